$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 10.379369
$ws.Range("H2").Value = 31.138107
$ws.Range("I2").Value = 0.01614698522449884
$ws.Range("J2").Value = 0.01614698522449883
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.2113696666666667
$ws.Range("N2").Value = 0.634109
$ws.Range("O2").Value = 0.03795977003925348
$ws.Range("P2").Value = 0.03795977003925347
$ws.Range("Q2").Value = 2.193883765740333
$ws.Range("R2").Value = 19.744953891663
$ws.Range("S2").Value = 0.0006129358459491995
$ws.Range("T2").Value = 0.0006129358459491993

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 10.379369
$ws.Range("H3").Value = 31.138107
$ws.Range("I3").Value = 0.01614698522449884
$ws.Range("J3").Value = 0.01614698522449883
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.589504333333333
$ws.Range("N3").Value = 4.768513
$ws.Range("O3").Value = 0.2854582680725092
$ws.Range("P3").Value = 0.2854582680725091
$ws.Range("Q3").Value = 16.49805200276567
$ws.Range("R3").Value = 148.482468024891
$ws.Range("S3").Value = 0.004609290436777834
$ws.Range("T3").Value = 0.004609290436777831

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 10.379369
$ws.Range("H4").Value = 31.138107
$ws.Range("I4").Value = 0.01614698522449884
$ws.Range("J4").Value = 0.01614698522449883
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3.767380666666666
$ws.Range("N4").Value = 11.302142
$ws.Range("O4").Value = 0.6765819618882374
$ws.Range("P4").Value = 0.6765819618882374
$ws.Range("Q4").Value = 39.10303410279933
$ws.Range("R4").Value = 351.927306925194
$ws.Range("S4").Value = 0.0109247589417718
$ws.Range("T4").Value = 0.0109247589417718

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 604.0312093333333
$ws.Range("H5").Value = 1812.093628
$ws.Range("I5").Value = 0.9396797639857967
$ws.Range("J5").Value = 0.9396797639857967
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.2113696666666667
$ws.Range("N5").Value = 0.634109
$ws.Range("O5").Value = 0.03795977003925348
$ws.Range("P5").Value = 0.03795977003925347
$ws.Range("Q5").Value = 127.6738753730502
$ws.Range("R5").Value = 1149.064878357452
$ws.Range("S5").Value = 0.03567002775144083
$ws.Range("T5").Value = 0.03567002775144082

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 604.0312093333333
$ws.Range("H6").Value = 1812.093628
$ws.Range("I6").Value = 0.9396797639857967
$ws.Range("J6").Value = 0.9396797639857967
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.589504333333333
$ws.Range("N6").Value = 4.768513
$ws.Range("O6").Value = 0.2854582680725092
$ws.Range("P6").Value = 0.2854582680725091
$ws.Range("Q6").Value = 960.1102247039072
$ws.Range("R6").Value = 8640.992022335166
$ws.Range("S6").Value = 0.2682393579701697
$ws.Range("T6").Value = 0.2682393579701697

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 604.0312093333333
$ws.Range("H7").Value = 1812.093628
$ws.Range("I7").Value = 0.9396797639857967
$ws.Range("J7").Value = 0.9396797639857967
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.767380666666666
$ws.Range("N7").Value = 11.302142
$ws.Range("O7").Value = 0.6765819618882374
$ws.Range("P7").Value = 0.6765819618882374
$ws.Range("Q7").Value = 2275.615500105686
$ws.Range("R7").Value = 20480.53950095118
$ws.Range("S7").Value = 0.6357703782641861
$ws.Range("T7").Value = 0.6357703782641861

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 28.39480333333333
$ws.Range("H8").Value = 85.18441
$ws.Range("I8").Value = 0.04417325078970442
$ws.Range("J8").Value = 0.04417325078970442
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.2113696666666667
$ws.Range("N8").Value = 0.634109
$ws.Range("O8").Value = 0.03795977003925348
$ws.Range("P8").Value = 0.03795977003925347
$ws.Range("Q8").Value = 6.001800115632222
$ws.Range("R8").Value = 54.01620104069
$ws.Range("S8").Value = 0.001676806441863452
$ws.Range("T8").Value = 0.001676806441863452

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 28.39480333333333
$ws.Range("H9").Value = 85.18441
$ws.Range("I9").Value = 0.04417325078970442
$ws.Range("J9").Value = 0.04417325078970442
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.589504333333333
$ws.Range("N9").Value = 4.768513
$ws.Range("O9").Value = 0.2854582680725092
$ws.Range("P9").Value = 0.2854582680725091
$ws.Range("Q9").Value = 45.13366294248111
$ws.Range("R9").Value = 406.2029664823301
$ws.Range("S9").Value = 0.01260961966556162
$ws.Range("T9").Value = 0.01260961966556162

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 28.39480333333333
$ws.Range("H10").Value = 85.18441
$ws.Range("I10").Value = 0.04417325078970442
$ws.Range("J10").Value = 0.04417325078970442
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.767380666666666
$ws.Range("N10").Value = 11.302142
$ws.Range("O10").Value = 0.6765819618882374
$ws.Range("P10").Value = 0.6765819618882374
$ws.Range("Q10").Value = 106.9740331118022
$ws.Range("R10").Value = 962.76629800622
$ws.Range("S10").Value = 0.02988682468227935
$ws.Range("T10").Value = 0.02988682468227935
